$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add "California" label in B1
$ws.Range("B1").Value = "California"

# Update the source/last-updated date in C1
$ws.Range("C1").Value = (Get-Date -Year 2022 -Month 1 -Day 4 -Hour 0 -Minute 0 -Second 0)
